# Regenerate handback-status report: swap the two e2e test-fixture file
# identifiers (and their recorded timestamps) that this report tracks.
#   0755463a-7804-4d15-8686-2fe2c663ed24  ->  7b1e7c19-20d0-4004-8d1a-136c09bbf563
#   dffb5ce4-5a26-49fb-b42e-b6af9bc09fba  ->  ffff72072017-1d0a-4d03-a789-57079cd6f6bb
# and the zh-cn / de-de generated .xlf artifact names + handoff/handback times.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "0755463a-7804-4d15-8686-2fe2c663ed24"
$newGuid1 = "7b1e7c19-20d0-4004-8d1a-136c09bbf563"
$oldGuid2 = "dffb5ce4-5a26-49fb-b42e-b6af9bc09fba"
$newGuid2 = "ffff72072017-1d0a-4d03-a789-57079cd6f6bb"

$newMd1 = "$newGuid1.md"
$newMd2 = "$newGuid2.md"

$newZhXlf = "$newGuid1.7bd3990a85e2cc549418b1a94481ffd36c08d550.zh-cn.xlf"
$newDeXlf = "$newGuid1.7bd3990a85e2cc549418b1a94481ffd36c08d550.de-de.xlf"

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = $newMd1
$ws.Range("B2").Value = "e2e\$newMd1"
$ws.Range("G2").Value = "2016-08-30 15:21:57"
$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = "e2e\$newMd2"
$ws.Range("G3").Value = "2016-08-30 15:21:57"

foreach ($h in $ws.Hyperlinks) {
    if ($h.TextToDisplay -eq "e2e\$oldGuid1.md") {
        $h.TextToDisplay = "e2e\$newMd1"
    } elseif ($h.TextToDisplay -eq "e2e\$oldGuid2.md") {
        $h.TextToDisplay = "e2e\$newMd2"
    }
}

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = $newMd1
$ws.Range("G2").Value = $newZhXlf
$ws.Range("H2").Value = "2016-08-30 15:21:52"
$ws.Range("I2").Value = $newMd1
$ws.Range("J2").Value = $newZhXlf
$ws.Range("K2").Value = "2016-08-30 15:22:26"

$ws.Range("A3").Value = $newMd2
$ws.Range("G3").Value = $newZhXlf
$ws.Range("H3").Value = "2016-08-30 15:21:52"
$ws.Range("I3").Value = $newMd2
$ws.Range("J3").Value = $newZhXlf
$ws.Range("K3").Value = "2016-08-30 15:22:26"

foreach ($h in $ws.Hyperlinks) {
    if ($h.TextToDisplay -eq "$oldGuid1.md") {
        $h.TextToDisplay = $newMd1
    } elseif ($h.TextToDisplay -eq "$oldGuid2.md") {
        $h.TextToDisplay = $newMd2
    }
}

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = $newMd1
$ws.Range("G2").Value = $newDeXlf
$ws.Range("H2").Value = "2016-08-30 15:21:57"
$ws.Range("I2").Value = $newMd1
$ws.Range("J2").Value = $newDeXlf
$ws.Range("K2").Value = "2016-08-30 15:22:33"

$ws.Range("A3").Value = $newMd2
$ws.Range("G3").Value = $newDeXlf
$ws.Range("H3").Value = "2016-08-30 15:21:57"
$ws.Range("I3").Value = $newMd2
$ws.Range("J3").Value = $newDeXlf
$ws.Range("K3").Value = "2016-08-30 15:22:33"

foreach ($h in $ws.Hyperlinks) {
    if ($h.TextToDisplay -eq "$oldGuid1.md") {
        $h.TextToDisplay = $newMd1
    } elseif ($h.TextToDisplay -eq "$oldGuid2.md") {
        $h.TextToDisplay = $newMd2
    }
}
